# Adds the new Bay row (row 47) that appears in the updated export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A47").Value = 1
$ws.Range("B47").Value = 4
$ws.Range("C47").Value = "fssd"
$ws.Range("D47").Value = "sdf"
$ws.Range("E47").Value = 1

# F47 mirrors the other blank "URLs" cells in this sheet, which are stored
# as empty text (not a truly blank/empty cell). A leading apostrophe forces
# the value to be stored as an empty string of type Text instead of
# clearing the cell; resetting the style afterwards drops the transient
# "quote prefix" formatting so the cell matches its neighbours exactly.
$ws.Range("F47").Value = "'"
$ws.Range("F47").Style = "Normal"
